$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert a new row above the current row 11 ("Description"), pushing
# Description/Purpose/Copyright/Immutable down by one row.
$ws1.Rows.Item(11).Insert()

# Copy the formatting (border/alignment/wrap) of the row above (Contact)
# onto the two new cells so they match the sheet's normal style instead of
# picking up a blank default style.
$ws1.Range("A10:B10").Copy()
$ws1.Range("A11:B11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New "Jurisdiction" metadata row (value left blank).
$ws1.Cells.Item(11, 1).Value = "Jurisdiction"
$ws1.Cells.Item(11, 2).Value = ""

# Refresh the "Date" metadata row with the new generation timestamp.
$ws1.Cells.Item(8, 2).Value = "2024-07-01T07:50:29+00:00"
